# Update duration system with new base values and multipliers.
#
# Base durations (Standard preset) change on the "Menu Mock" sheet (column D,
# the current/default value of each OptionDuration row), and every row's
# allowed-values list (column E) plus the shared provider list on the
# "Providers" sheet (B11, CustomDurationProvider) gains three new entries
# (2.75s, 3.25s, 3.5s) inserted in sorted order.

$wb = $excel.ActiveWorkbook
$wsMenu = $wb.Worksheets.Item("Menu Mock")
$wsProviders = $wb.Worksheets.Item("Providers")

$newAllowedValues = "0.5s | 0.6s | 0.72s | 0.75s | 0.90s | 1.0s | 1.125s | 1.2s | 1.25s | 1.4s | 1.5s | 1.68s | 1.8s | 1.875s | 2.0s | 2.1s | 2.16s | 2.25s | 2.4s | 2.5s | 2.7s | 2.75s | 2.8s | 3.0s | 3.25s | 3.5s | 3.6s | 3.75s | 4.0s | 4.2s | 4.5s | 5.0s | 5.4s | 6.0s | 6.25s | 7.0s | 7.5s | 9.0s | 10.0s"

# Row -> new base duration value, per the commit message:
#   Row 44 CategoryCustomBasic (Basic Kill)       -> 2.5s
#   Row 53 CategoryCustomCritical (Critical)      -> 3.0s
#   Row 62 CategoryCustomDismemberment            -> 2.0s
#   Row 71 CategoryCustomDecapitation             -> 3.25s
#   Row 80 CategoryCustomLastEnemy (Last Enemy)   -> 2.75s
#   Row 88 CategoryCustomLastStand (Last Stand)   -> 3.5s
#   Row 96 CategoryCustomParry (Parry)            -> 1.5s
$durationRows = @{
    44 = "2.5s"
    53 = "3.0s"
    62 = "2.0s"
    71 = "3.25s"
    80 = "2.75s"
    88 = "3.5s"
    96 = "1.5s"
}

foreach ($row in $durationRows.Keys) {
    $wsMenu.Range("D$row").Value = $durationRows[$row]
    $wsMenu.Range("E$row").Value = $newAllowedValues
}

# Update the shared CustomDurationProvider values list on the Providers sheet.
$wsProviders.Range("B11").Value = $newAllowedValues
